# Entering storage info and day 13 cohort 2018-08-28 treatments
$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------------
# 1. Duplicate the "R4_P3" sheet (the last R4 sheet) to create a new sheet
#    for the day 13 cohort, placed right after it, and rename it "R4_P4".
# ---------------------------------------------------------------------------
$oldSheet = $wb.Worksheets.Item("R4_P3")
$oldSheet.Copy([System.Reflection.Missing]::Value, $oldSheet)
$newSheet = $wb.Worksheets.Item("R4_P3 (2)")
$newSheet.Name = "R4_P4"

# ---------------------------------------------------------------------------
# 2. Fill in the storage-info header fields on the new "R4_P4" sheet with
#    the day 13 cohort information.
# ---------------------------------------------------------------------------
$newSheet.Range("B1").Value = "2018-09-12 started"
$newSheet.Range("B2").Value = "Genetic control for apple race, 2018 collection year, day 13(2018-08-28)"
$newSheet.Range("B4").Value = "4"
$newSheet.Range("B5").Clear()
$newSheet.Range("B6").Value = "nb#004, page 162,158"

# ---------------------------------------------------------------------------
# 3. Clear the old day-10-cohort tube grid and highlight the still-empty
#    "fill in later" rows of the new cohort in yellow.
# ---------------------------------------------------------------------------
$newSheet.Range("C9:K9").ClearContents()
$newSheet.Range("C9:K9").Interior.Color = 65535
$newSheet.Range("D9").Value = "tubes for cohortday 12 genetic control; fill in later"

$newSheet.Range("C10:I10").ClearContents()
$newSheet.Range("C10:I10").Interior.Color = 65535
$newSheet.Range("J10").ClearContents()

# ---------------------------------------------------------------------------
# 4. Enter the day 13 cohort tube labels (A13-1 .. A13-30) into the grid.
# ---------------------------------------------------------------------------
$newSheet.Range("C11").Value = "A13-1"
$newSheet.Range("D11").Value = "A13-2"
$newSheet.Range("E11").Value = "A13-3"
$newSheet.Range("F11").Value = "A13-4"
$newSheet.Range("G11").Value = "A13-5"
$newSheet.Range("H11").Value = "A13-6"
$newSheet.Range("I11").Value = "A13-7"
$newSheet.Range("J11").Value = "A13-8"
$newSheet.Range("K11").Value = "A13-9"

$newSheet.Range("C12").Value = "A13-10"
$newSheet.Range("D12").Value = "A13-11"
$newSheet.Range("E12").Value = "A13-12"
$newSheet.Range("F12").Value = "A13-13"
$newSheet.Range("G12").Value = "A13-14"
$newSheet.Range("H12").Value = "A13-15"
$newSheet.Range("I12").Value = "A13-16"
$newSheet.Range("J12").Value = "A13-17"
$newSheet.Range("K12").Value = "A13-18"

$newSheet.Range("C13").Value = "A13-19"
$newSheet.Range("D13").Value = "A13-20"
$newSheet.Range("E13").Value = "A13-21"
$newSheet.Range("F13").Value = "A13-22"
$newSheet.Range("G13").Value = "A13-23"
$newSheet.Range("H13").Value = "A13-24"
$newSheet.Range("I13").Value = "A13-25"
$newSheet.Range("J13").Value = "A13-26"
$newSheet.Range("K13").Value = "A13-27"

$newSheet.Range("C14").Value = "A13-28"
$newSheet.Range("D14").Value = "A13-29"
$newSheet.Range("E14").Value = "A13-30"
$newSheet.Range("F14:K14").ClearContents()

$newSheet.Range("C15:K15").ClearContents()
$newSheet.Range("C16:G16").ClearContents()

# ---------------------------------------------------------------------------
# 5. View/selection tweaks: the new sheet becomes the active tab with B5
#    selected; the old "R4_P3" sheet is left with its whole grid selected.
# ---------------------------------------------------------------------------
$newSheet.Activate()
$newSheet.Range("B5").Select()

$oldSheet.Activate()
$oldSheet.Cells.Select()

$newSheet.Activate()
